$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters to indices: B=2 C=3 D=4 E=5
# D column holds price values that may look numeric; force text storage
# by setting NumberFormat to "@" before assignment, then restoring the
# original cell style afterward so the saved style index is unchanged.

# Row 2
$cell = $ws.Cells.Item(2, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "22.384.33"
$cell.Style = $origStyle
$ws.Cells.Item(2, 5).Value = "  +9.09%  "

# Row 3
$cell = $ws.Cells.Item(3, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.601.37"
$cell.Style = $origStyle
$ws.Cells.Item(3, 5).Value = "  +8.66%  "

# Row 4
$cell = $ws.Cells.Item(4, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.007"
$cell.Style = $origStyle
$ws.Cells.Item(4, 5).Value = "  -0.21%  "

# Row 5
$cell = $ws.Cells.Item(5, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "305.59"
$cell.Style = $origStyle
$ws.Cells.Item(5, 5).Value = "  +10.15%  "

# Row 6
$cell = $ws.Cells.Item(6, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9934"
$cell.Style = $origStyle
$ws.Cells.Item(6, 5).Value = "  +4.15%  "

# Row 7
$cell = $ws.Cells.Item(7, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.3661"
$cell.Style = $origStyle
$ws.Cells.Item(7, 5).Value = "  +1.46%  "

# Row 8
$cell = $ws.Cells.Item(8, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.3392"
$cell.Style = $origStyle
$ws.Cells.Item(8, 5).Value = "  +11.02%  "

# Row 9
$cell = $ws.Cells.Item(9, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "42.22"
$cell.Style = $origStyle
$ws.Cells.Item(9, 5).Value = "  +7.23%  "

# Row 10
$cell = $ws.Cells.Item(10, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.127"
$cell.Style = $origStyle
$ws.Cells.Item(10, 5).Value = "  +6.68%  "

# Row 11
$cell = $ws.Cells.Item(11, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.07034"
$cell.Style = $origStyle
$ws.Cells.Item(11, 5).Value = "  +5.98%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -0.03%  "

# Row 13
$cell = $ws.Cells.Item(13, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "19.67"
$cell.Style = $origStyle
$ws.Cells.Item(13, 5).Value = "  +8.76%  "

# Row 14
$cell = $ws.Cells.Item(14, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.899"
$cell.Style = $origStyle
$ws.Cells.Item(14, 5).Value = "  +7.15%  "

# Row 15
$cell = $ws.Cells.Item(15, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.608"
$cell.Style = $origStyle
$ws.Cells.Item(15, 5).Value = "  +7.03%  "

# Row 16
$cell = $ws.Cells.Item(16, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.605.67"
$cell.Style = $origStyle
$ws.Cells.Item(16, 5).Value = "  +8.85%  "

# Row 17
$cell = $ws.Cells.Item(17, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.00001079"
$cell.Style = $origStyle
$ws.Cells.Item(17, 5).Value = "  +5.05%  "

# Row 18
$cell = $ws.Cells.Item(18, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9942"
$cell.Style = $origStyle
$ws.Cells.Item(18, 5).Value = "  +3.94%  "

# Row 19
$cell = $ws.Cells.Item(19, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06634"
$cell.Style = $origStyle
$ws.Cells.Item(19, 5).Value = "  +11.60%  "

# Row 20
$cell = $ws.Cells.Item(20, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "77.58"
$cell.Style = $origStyle
$ws.Cells.Item(20, 5).Value = "  +12.35%  "

# Row 21
$cell = $ws.Cells.Item(21, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.000"
$cell.Style = $origStyle
$ws.Cells.Item(21, 5).Value = "  +9.26%  "

# Row 22
$cell = $ws.Cells.Item(22, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "15.96"
$cell.Style = $origStyle
$ws.Cells.Item(22, 5).Value = "  +10.33%  "

# Row 23
$cell = $ws.Cells.Item(23, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.81"
$cell.Style = $origStyle
$ws.Cells.Item(23, 5).Value = "  +6.06%  "

# Row 24
$cell = $ws.Cells.Item(24, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "22.449.21"
$cell.Style = $origStyle
$ws.Cells.Item(24, 5).Value = "  +9.17%  "

# Row 25
$cell = $ws.Cells.Item(25, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.386"
$cell.Style = $origStyle
$ws.Cells.Item(25, 5).Value = "  +6.00%  "

# Row 26
$cell = $ws.Cells.Item(26, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.547"
$cell.Style = $origStyle
$ws.Cells.Item(26, 5).Value = "  +19.64%  "

# Row 27
$cell = $ws.Cells.Item(27, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "149.37"
$cell.Style = $origStyle
$ws.Cells.Item(27, 5).Value = "  +4.35%  "

# Row 28
$cell = $ws.Cells.Item(28, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "19.36"
$cell.Style = $origStyle
$ws.Cells.Item(28, 5).Value = "  +12.90%  "

# Row 29
$cell = $ws.Cells.Item(29, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.786.67"
$cell.Style = $origStyle
$ws.Cells.Item(29, 5).Value = "  +9.16%  "

# Row 30
$cell = $ws.Cells.Item(30, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "122.41"
$cell.Style = $origStyle
$ws.Cells.Item(30, 5).Value = "  +7.84%  "

# Row 31
$cell = $ws.Cells.Item(31, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.051"
$cell.Style = $origStyle
$ws.Cells.Item(31, 5).Value = "  +2.80%  "

# Row 32
$cell = $ws.Cells.Item(32, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.110"
$cell.Style = $origStyle
$ws.Cells.Item(32, 5).Value = "  +22.24%  "

# Row 33
$cell = $ws.Cells.Item(33, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9431"
$cell.Style = $origStyle
$ws.Cells.Item(33, 5).Value = "  +16.97%  "

# Row 34
$cell = $ws.Cells.Item(34, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.690"
$cell.Style = $origStyle
$ws.Cells.Item(34, 5).Value = "  +11.69%  "

# Row 35
$cell = $ws.Cells.Item(35, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.08207"
$cell.Style = $origStyle
$ws.Cells.Item(35, 5).Value = "  +2.79%  "

# Row 36
$cell = $ws.Cells.Item(36, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.88"
$cell.Style = $origStyle
$ws.Cells.Item(36, 5).Value = "  +15.08%  "

# Row 37
$cell = $ws.Cells.Item(37, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.205"
$cell.Style = $origStyle
$ws.Cells.Item(37, 5).Value = "  +10.52%  "

# Row 38
$cell = $ws.Cells.Item(38, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.256"
$cell.Style = $origStyle
$ws.Cells.Item(38, 5).Value = "  +3.36%  "

# Row 39
$cell = $ws.Cells.Item(39, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "8.516"
$cell.Style = $origStyle
$ws.Cells.Item(39, 5).Value = "  +14.62%  "

# Row 40
$cell = $ws.Cells.Item(40, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06052"
$cell.Style = $origStyle
$ws.Cells.Item(40, 5).Value = "  +3.72%  "

# Row 41
$cell = $ws.Cells.Item(41, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.02204"
$cell.Style = $origStyle
$ws.Cells.Item(41, 5).Value = "  +7.55%  "

# Row 42
$cell = $ws.Cells.Item(42, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.2018"
$cell.Style = $origStyle

# Row 43
$cell = $ws.Cells.Item(43, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9936"
$cell.Style = $origStyle
$ws.Cells.Item(43, 5).Value = "  +3.87%  "

# Row 44
$cell = $ws.Cells.Item(44, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.5885"
$cell.Style = $origStyle
$ws.Cells.Item(44, 5).Value = "  +11.46%  "

# Row 45
$cell = $ws.Cells.Item(45, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.841"
$cell.Style = $origStyle
$ws.Cells.Item(45, 5).Value = "  +9.09%  "

# Row 46
$cell = $ws.Cells.Item(46, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "13.05"
$cell.Style = $origStyle
$ws.Cells.Item(46, 5).Value = "  +7.18%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Quant"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$cell = $ws.Cells.Item(47, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "127.04"
$cell.Style = $origStyle
$ws.Cells.Item(47, 5).Value = "  +7.48%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "Decentraland"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$cell = $ws.Cells.Item(48, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.5664"
$cell.Style = $origStyle
$ws.Cells.Item(48, 5).Value = "  +9.19%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +8.20%  "

# Row 50
$cell = $ws.Cells.Item(50, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.06803"
$cell.Style = $origStyle
$ws.Cells.Item(50, 5).Value = "  +5.11%  "

# Row 51
$cell = $ws.Cells.Item(51, 4)
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "73.39"
$cell.Style = $origStyle
$ws.Cells.Item(51, 5).Value = "  +8.96%  "
